$d = $word.ActiveDocument

# 1) Remove the ".3" suffix run after "Porta LEDs" (merge back into a single run)
$d.Content.Find.Execute("Porta LEDs.3", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Porta LEDs", 2)

# 2) Prefix "Jack DC hembra para PCB" with a new run containing "1 x "
$rng = $d.Content
$rng.Find.Execute("Jack DC hembra para PCB", $true, $false, $false, $false, $false,
                   $true, 1, $false, "", 0)
$ins = $rng.Duplicate
$ins.Collapse(1)
$ins.Text = "1 x "
$ins.Bold = 1
$ins.Bold = 0

# 3) Prefix "Jack RJ45" with a new run containing "1 x "
$rng2 = $d.Content
$rng2.Find.Execute("Jack RJ45", $true, $false, $false, $false, $false,
                    $true, 1, $false, "", 0)
$ins2 = $rng2.Duplicate
$ins2.Collapse(1)
$ins2.Text = "1 x "
$ins2.Bold = 1
$ins2.Bold = 0

# 4) Simplify "1 x SPST o 1 x de enclavamiento" to a single run "1 x SPST"
$d.Content.Find.Execute("1 x SPST o 1 x de enclavamiento", $true, $false, $false, $false, $false,
                         $true, 1, $false, "1 x SPST", 2)
